$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.038.54"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.582.08"
$ws.Range("E3").Value = "  +8.79%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'305.70"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "'99.60"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("D7").Value = "'0.595"
$ws.Range("E7").Value = "  +5.54%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.573"
$ws.Range("E9").Value = "  +13.04%  "
$ws.Range("D10").Value = "'38.45"
$ws.Range("E10").Value = "  +12.05%  "
$ws.Range("E11").Value = "  +6.74%  "
$ws.Range("D12").Value = "'8.30"
$ws.Range("E12").Value = "  +16.95%  "
$ws.Range("D13").Value = "2.978.33"
$ws.Range("E13").Value = "  +8.80%  "
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "2.604.31"
$ws.Range("E15").Value = "  +9.56%  "
$ws.Range("D16").Value = "'0.900"
$ws.Range("E16").Value = "  +9.54%  "
$ws.Range("D17").Value = "'14.80"
$ws.Range("E17").Value = "  +8.01%  "
$ws.Range("D18").Value = "46.174.46"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  +6.77%  "
$ws.Range("D20").Value = "'12.92"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").Value = "'6.62"
$ws.Range("E21").Value = "  +10.09%  "
$ws.Range("D22").Value = "'71.03"
$ws.Range("E22").Value = "  +6.38%  "
$ws.Range("D23").Value = "'253.06"
$ws.Range("E23").Value = "  +4.17%  "
$ws.Range("E24").Value = "  +7.42%  "
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  +13.78%  "
$ws.Range("D26").Value = "'28.23"
$ws.Range("E26").Value = "  +34.57%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'10.42"
$ws.Range("E28").Value = "  +7.39%  "
$ws.Range("D29").Value = "'39.56"
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").Value = "'2.26"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").Value = "'6.04"
$ws.Range("E31").Value = "  +9.44%  "
$ws.Range("D32").Value = "'3.67"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").Value = "'2.92"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("E34").Value = "  +19.52%  "
$ws.Range("D35").Value = "'152.60"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("D36").Value = "'0.0822"
$ws.Range("E36").Value = "  +7.13%  "
$ws.Range("E37").Value = "  +3.34%  "
$ws.Range("E38").Value = "  +5.39%  "
$ws.Range("D39").Value = "'16.08"
$ws.Range("E39").Value = "  +8.13%  "
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("D41").Value = "'3.59"
$ws.Range("E41").Value = "  +12.52%  "
$ws.Range("D42").Value = "'0.0319"
$ws.Range("E42").Value = "  +7.34%  "
$ws.Range("D43").Value = "2.051.91"
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("D44").Value = "'19.83"
$ws.Range("E44").Value = "  +40.33%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "'90.76"
$ws.Range("E46").Value = "  -4.08%  "
$ws.Range("E47").Value = "  +9.35%  "
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "'108.19"
$ws.Range("E49").Value = "  +9.45%  "
$ws.Range("E50").Value = "  +8.36%  "
$ws.Range("D51").Value = "2.836.34"
$ws.Range("E51").Value = "  +8.77%  "
